$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(24, 42602.513599537036, "Noun", 2770, 75, 8, 2, 1, 66, 33, 4, 3, 57, 42),
    @(25, 42602.516712962963, "Noun", 3058, 75, 8, 2, 1, 66, 33, 4, 3, 57, 42),
    @(26, 42602.524062500001, "Noun", 3069, 75, 8, 2, 1, 66, 33, 4, 3, 57, 42)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
    $ws.Cells.Item($r, 9).Value = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
    $ws.Cells.Item($r, 12).Value = $row[12]
    $ws.Cells.Item($r, 13).Value = $row[13]
}
